# Replace the two-digit multiplication problems/answers in the table cells
# with the values from the updated worksheet, per the commit diff.
# Each old string is unique in the document, so a straightforward
# Find/Replace (MatchCase, not whole-word since strings contain symbols
# like "×" and "=") is sufficient and safe.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "60×62=3720" "13×70=910"
Replace-Text "22×30=660" "73×52=3796"
Replace-Text "67×74=4958" "51×66=3366"
Replace-Text "50×58=2900" "14×58=812"
Replace-Text "92×52=4784" "27×15=405"
Replace-Text "73×98=7154" "13×27=351"
Replace-Text "55×73=4015" "97×79=7663"
Replace-Text "37×93=3441" "29×49=1421"
Replace-Text "61×68=4148" "49×34=1666"
Replace-Text "82×54=4428" "28×31=868"
Replace-Text "17×31=527" "83×59=4897"
Replace-Text "59×51=3009" "85×91=7735"
Replace-Text "84×55=4620" "42×98=4116"
Replace-Text "45×94=4230" "26×98=2548"
Replace-Text "35×83=2905" "45×74=3330"
Replace-Text "82×37=3034" "15×87=1305"
Replace-Text "32×78=2496" "92×97=8924"
Replace-Text "60×98=5880" "97×38=3686"
Replace-Text "14×67=938" "59×60=3540"
Replace-Text "98×43=4214" "69×64=4416"
Replace-Text "55×97=5335" "20×87=1740"
Replace-Text "72×43=3096" "24×33=792"
Replace-Text "87×13=1131" "27×49=1323"
Replace-Text "89×70=6230" "74×32=2368"
Replace-Text "71×93=6603" "60×29=1740"
